$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewPolicy")

# Fix up TestColCount (T2) for the existing "newpol_001" row, then append a new
# "newpol_002" policy-creation test row (row 3) mirroring row 2's layout.
$ws.Cells.Item(2, 20).Value = 50

$ws.Cells.Item(3, 1).Value  = "newpol_002"
$ws.Cells.Item(3, 2).Value  = "Create New General Policy"
$ws.Cells.Item(3, 3).Value  = "PolicyData"
$ws.Cells.Item(3, 4).Value  = "NewPolicy"
$ws.Cells.Item(3, 10).Value = "ResNewPolicy"
$ws.Cells.Item(3, 11).Value = "NewPolicy"
$ws.Cells.Item(3, 12).Value = 4
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 3
$ws.Cells.Item(3, 15).Value = "ResNewPolicy"
$ws.Cells.Item(3, 16).Value = "NewPolicy"
$ws.Cells.Item(3, 17).Value = 4
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 1
$ws.Cells.Item(3, 20).Value = 50
$ws.Cells.Item(3, 21).Value = 1

# Bring the new row into view / selection, matching the post-edit UI state.
$ws.Activate()
[void]$ws.Range("A3").Select()
